# Auto-generated Excel COM-interop script to apply the cell value changes
# described by the commit diff for Cactuar_Profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 350.16666
$ws.Range("I33").Value = 374.75
$ws.Range("J33").Value = 301
$ws.Range("K33").Value = 374.75
$ws.Range("L33").Value = 301
$ws.Range("M33").Value = -145.75
$ws.Range("N33").Value = -759

$ws.Range("H40").Value = 21926.834
$ws.Range("I40").Value = 23659
$ws.Range("J40").Value = 19761.625
$ws.Range("K40").Value = 23659
$ws.Range("L40").Value = 19761.625
$ws.Range("M40").Value = -23484
$ws.Range("N40").Value = -20111.625

$ws.Range("H51").Value = 7564.788
$ws.Range("J51").Value = 8199.538
$ws.Range("L51").Value = 8199.538
$ws.Range("N51").Value = -9167.538

$ws.Range("H95").Value = 30312
$ws.Range("J95").Value = 30312
$ws.Range("L95").Value = 30312
$ws.Range("N95").Value = -35804

$ws.Range("H98").Value = 1757.091
$ws.Range("I98").Value = 1757.091
$ws.Range("K98").Value = 1757.091
$ws.Range("M98").Value = -259.0909999999999

$ws.Range("H122").Value = 1757.091
$ws.Range("I122").Value = 1757.091
$ws.Range("K122").Value = 5271.272999999999
$ws.Range("M122").Value = -2821.272999999999

$ws.Range("H132").Value = 170941.6
$ws.Range("I132").Value = 254522.61
$ws.Range("K132").Value = 763567.83
$ws.Range("M132").Value = -761037.83

$ws.Range("H135").Value = 3507.9565
$ws.Range("I135").Value = 1892.35
$ws.Range("J135").Value = 14278.667
$ws.Range("K135").Value = 17031.15
$ws.Range("L135").Value = 128508.003
$ws.Range("M135").Value = -14496.15
$ws.Range("N135").Value = -133578.003

$ws.Range("H137").Value = 13736105
$ws.Range("I137").Value = 529259.9399999999
$ws.Range("J137").Value = 55557784
$ws.Range("K137").Value = 1587779.82
$ws.Range("L137").Value = 166673352
$ws.Range("M137").Value = -1585229.82
$ws.Range("N137").Value = -166678452

$ws.Range("H138").Value = 5744.6562
$ws.Range("I138").Value = 2294.889
$ws.Range("J138").Value = 7094.5654
$ws.Range("K138").Value = 6884.667
$ws.Range("L138").Value = 21283.6962
$ws.Range("M138").Value = -1744.667
$ws.Range("N138").Value = -31563.6962

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10562.3955
$ws.Range("I32").Value = 10304.09
$ws.Range("K32").Value = 10304.09
$ws.Range("M32").Value = -10017.09

$ws.Range("H45").Value = 3832.1538
$ws.Range("I45").Value = 1974.1428
$ws.Range("J45").Value = 5999.8335
$ws.Range("K45").Value = 1974.1428
$ws.Range("L45").Value = 5999.8335
$ws.Range("M45").Value = -1597.1428
$ws.Range("N45").Value = -6753.8335

$ws.Range("H61").Value = 3322.3784
$ws.Range("I61").Value = 2701.074
$ws.Range("J61").Value = 4999.9
$ws.Range("K61").Value = 2701.074
$ws.Range("L61").Value = 4999.9
$ws.Range("M61").Value = -2489.074
$ws.Range("N61").Value = -5423.9

$ws.Range("H74").Value = 1398.3422
$ws.Range("I74").Value = 1038.0385
$ws.Range("J74").Value = 2179
$ws.Range("K74").Value = 1038.0385
$ws.Range("L74").Value = 2179
$ws.Range("M74").Value = -164.0385000000001
$ws.Range("N74").Value = -3927

$ws.Range("H77").Value = 1398.3422
$ws.Range("I77").Value = 1038.0385
$ws.Range("J77").Value = 2179
$ws.Range("K77").Value = 5190.192500000001
$ws.Range("L77").Value = 10895
$ws.Range("M77").Value = -822.192500000001
$ws.Range("N77").Value = -19631

$ws.Range("H95").Value = 21421.4
$ws.Range("I95").Value = 10000
$ws.Range("J95").Value = 24276.75
$ws.Range("K95").Value = 10000
$ws.Range("L95").Value = 24276.75
$ws.Range("M95").Value = -7254
$ws.Range("N95").Value = -29768.75

$ws.Range("H132").Value = 12945.35
$ws.Range("I132").Value = 14158.637
$ws.Range("J132").Value = 9608.8125
$ws.Range("K132").Value = 42475.911
$ws.Range("L132").Value = 28826.4375
$ws.Range("M132").Value = -39945.911
$ws.Range("N132").Value = -33886.4375

$ws.Range("H136").Value = 3322.3784
$ws.Range("I136").Value = 2701.074
$ws.Range("J136").Value = 4999.9
$ws.Range("K136").Value = 8103.222
$ws.Range("L136").Value = 14999.7
$ws.Range("M136").Value = -5553.222
$ws.Range("N136").Value = -20099.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6557.24
$ws.Range("I20").Value = 6303
$ws.Range("J20").Value = 7009.222
$ws.Range("K20").Value = 6303
$ws.Range("L20").Value = 7009.222
$ws.Range("M20").Value = -6056
$ws.Range("N20").Value = -7503.222

$ws.Range("H22").Value = 600
$ws.Range("I22").Value = 600
$ws.Range("K22").Value = 600
$ws.Range("M22").Value = -427

$ws.Range("H86").Value = 1068.75
$ws.Range("I86").Value = 1078.5714
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 1078.5714
$ws.Range("L86").Value = 1000
$ws.Range("M86").Value = 44.42859999999996
$ws.Range("N86").Value = -3246

$ws.Range("H89").Value = 1068.75
$ws.Range("I89").Value = 1078.5714
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 5392.857
$ws.Range("L89").Value = 5000
$ws.Range("M89").Value = 223.143
$ws.Range("N89").Value = -16232

$ws.Range("H134").Value = 3083.5625
$ws.Range("I134").Value = 2666.9285
$ws.Range("K134").Value = 8000.7855
$ws.Range("M134").Value = -5465.7855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1276.2
$ws.Range("J22").Value = 1611
$ws.Range("L22").Value = 1611
$ws.Range("N22").Value = -2311

$ws.Range("H31").Value = 17548012
$ws.Range("I31").Value = 40002480
$ws.Range("J31").Value = 5457.9375
$ws.Range("K31").Value = 40002480
$ws.Range("L31").Value = 5457.9375
$ws.Range("M31").Value = -40002185
$ws.Range("N31").Value = -6047.9375

$ws.Range("H32").Value = 12501499
$ws.Range("I32").Value = 12501499
$ws.Range("K32").Value = 12501499
$ws.Range("M32").Value = -12501183

$ws.Range("H34").Value = 17548012
$ws.Range("I34").Value = 40002480
$ws.Range("J34").Value = 5457.9375
$ws.Range("K34").Value = 40002480
$ws.Range("L34").Value = 5457.9375
$ws.Range("M34").Value = -40002278
$ws.Range("N34").Value = -5861.9375

$ws.Range("H58").Value = 557905.3
$ws.Range("I58").Value = 2618.2856
$ws.Range("J58").Value = 2501410
$ws.Range("K58").Value = 2618.2856
$ws.Range("L58").Value = 2501410
$ws.Range("M58").Value = -2415.2856
$ws.Range("N58").Value = -2501816

$ws.Range("H62").Value = 43470.5
$ws.Range("I62").Value = 20784.166
$ws.Range("J62").Value = 77500
$ws.Range("K62").Value = 20784.166
$ws.Range("L62").Value = 77500
$ws.Range("M62").Value = -20160.166
$ws.Range("N62").Value = -78748

$ws.Range("H65").Value = 43470.5
$ws.Range("I65").Value = 20784.166
$ws.Range("J65").Value = 77500
$ws.Range("K65").Value = 103920.83
$ws.Range("L65").Value = 387500
$ws.Range("M65").Value = -100800.83
$ws.Range("N65").Value = -393740

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H94").Value = 2324.158
$ws.Range("J94").Value = 2634.0667
$ws.Range("L94").Value = 2634.0667
$ws.Range("N94").Value = -3536.0667

$ws.Range("H120").Value = 45991.332
$ws.Range("J120").Value = 45991.332
$ws.Range("L120").Value = 45991.332
$ws.Range("N120").Value = -53249.332

$ws.Range("H121").Value = 58316.332
$ws.Range("J121").Value = 58316.332
$ws.Range("L121").Value = 58316.332
$ws.Range("N121").Value = -60936.332

$ws.Range("H122").Value = 2891.6875
$ws.Range("I122").Value = 1465.5652
$ws.Range("K122").Value = 4396.6956
$ws.Range("M122").Value = -1946.6956

$ws.Range("H132").Value = 17792.46
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 17792.46
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 53377.38
$ws.Range("N132").Value = -58437.38
$ws.Range("M132").ClearContents()

$ws.Range("H134").Value = 1664.3846
$ws.Range("I134").Value = 1535.0227
$ws.Range("K134").Value = 4605.0681
$ws.Range("M134").Value = -2070.0681

$ws.Range("H136").Value = 557905.3
$ws.Range("I136").Value = 2618.2856
$ws.Range("J136").Value = 2501410
$ws.Range("K136").Value = 7854.8568
$ws.Range("L136").Value = 7504230
$ws.Range("M136").Value = -5304.8568
$ws.Range("N136").Value = -7509330

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 45.714287
$ws.Range("I16").Value = 43.6
$ws.Range("J16").Value = 51
$ws.Range("K16").Value = 130.8
$ws.Range("L16").Value = 153
$ws.Range("M16").Value = 42.19999999999999
$ws.Range("N16").Value = -499

$ws.Range("H44").Value = 16667683
$ws.Range("J44").Value = 1833.3334
$ws.Range("L44").Value = 5500.0002
$ws.Range("N44").Value = -6296.0002

$ws.Range("H47").Value = 5788.077
$ws.Range("I47").Value = 5780.125
$ws.Range("K47").Value = 17340.375
$ws.Range("M47").Value = -16909.375

$ws.Range("H136").Value = 2215
$ws.Range("J136").Value = 3500
$ws.Range("L136").Value = 10500
$ws.Range("N136").Value = -20700

$ws.Range("H137").Value = 141667200
$ws.Range("J137").Value = 50000600
$ws.Range("L137").Value = 150001800
$ws.Range("N137").Value = -150012000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 50000
$ws.Range("J52").Value = 50000
$ws.Range("L52").Value = 50000
$ws.Range("N52").Value = -50518

$ws.Range("H70").Value = 8029.6
$ws.Range("I70").Value = 8041.7856
$ws.Range("K70").Value = 8041.7856
$ws.Range("M70").Value = -7771.7856

$ws.Range("H73").Value = 8029.6
$ws.Range("I73").Value = 8041.7856
$ws.Range("K73").Value = 8041.7856
$ws.Range("M73").Value = -7105.7856

$ws.Range("H98").Value = 35333
$ws.Range("J98").Value = 35333
$ws.Range("L98").Value = 35333
$ws.Range("N98").Value = -41323

$ws.Range("H102").Value = 5647.1025
$ws.Range("I102").Value = 4880.7666
$ws.Range("K102").Value = 4880.7666
$ws.Range("M102").Value = -3258.7666

$ws.Range("H105").Value = 43667
$ws.Range("J105").Value = 43667
$ws.Range("L105").Value = 43667
$ws.Range("N105").Value = -50655

$ws.Range("H107").Value = 589.5294
$ws.Range("I107").Value = 394.5
$ws.Range("J107").Value = 1499.6666
$ws.Range("K107").Value = 394.5
$ws.Range("L107").Value = 1499.6666
$ws.Range("M107").Value = 1525.5
$ws.Range("N107").Value = -5339.6666

$ws.Range("H122").Value = 383423.7
$ws.Range("I122").Value = 581074.0600000001
$ws.Range("K122").Value = 1743222.18
$ws.Range("M122").Value = -1740772.18

$ws.Range("H127").Value = 80241.664
$ws.Range("J127").Value = 80241.664
$ws.Range("L127").Value = 80241.664
$ws.Range("N127").Value = -90161.664

$ws.Range("H132").Value = 406348.8
$ws.Range("I132").Value = 83217.96000000001
$ws.Range("J132").Value = 2506699.2
$ws.Range("K132").Value = 249653.88
$ws.Range("L132").Value = 7520097.600000001
$ws.Range("M132").Value = -247123.88
$ws.Range("N132").Value = -7525157.600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7019.9
$ws.Range("I7").Value = 4449.5
$ws.Range("J7").Value = 7662.5
$ws.Range("K7").Value = 4449.5
$ws.Range("L7").Value = 7662.5
$ws.Range("M7").Value = -4337.5
$ws.Range("N7").Value = -7886.5

$ws.Range("H16").Value = 2699.9583
$ws.Range("I16").Value = 2049.1333
$ws.Range("J16").Value = 3784.6667
$ws.Range("K16").Value = 2049.1333
$ws.Range("L16").Value = 3784.6667
$ws.Range("M16").Value = -1879.1333
$ws.Range("N16").Value = -4124.6667

$ws.Range("H22").Value = 1170.6666
$ws.Range("J22").Value = 1358.3334
$ws.Range("L22").Value = 1358.3334
$ws.Range("N22").Value = -1948.3334

$ws.Range("H27").Value = 1170.6666
$ws.Range("J27").Value = 1358.3334
$ws.Range("L27").Value = 1358.3334
$ws.Range("N27").Value = -1572.3334

$ws.Range("H40").Value = 8339679.5
$ws.Range("I40").Value = 31252250
$ws.Range("K40").Value = 31252250
$ws.Range("M40").Value = -31252114

$ws.Range("H47").Value = 59995
$ws.Range("J47").Value = 59995
$ws.Range("L47").Value = 59995
$ws.Range("N47").Value = -60975

$ws.Range("H52").Value = 59995
$ws.Range("J52").Value = 59995
$ws.Range("L52").Value = 59995
$ws.Range("N52").Value = -60461

$ws.Range("H54").Value = 36287.5
$ws.Range("J54").Value = 36287.5
$ws.Range("L54").Value = 36287.5
$ws.Range("N54").Value = -37575.5

$ws.Range("H100").Value = 2998.5
$ws.Range("I100").Value = 2998.5
$ws.Range("K100").Value = 2998.5
$ws.Range("M100").Value = -2457.5

$ws.Range("H122").Value = 8747.75
$ws.Range("I122").Value = 4818
$ws.Range("K122").Value = 14454
$ws.Range("M122").Value = -12004

$ws.Range("H126").Value = 7019.9
$ws.Range("I126").Value = 4449.5
$ws.Range("J126").Value = 7662.5
$ws.Range("K126").Value = 13348.5
$ws.Range("L126").Value = 22987.5
$ws.Range("M126").Value = -10878.5
$ws.Range("N126").Value = -27927.5

$ws.Range("H132").Value = 4392.4062
$ws.Range("I132").Value = 4745.619
$ws.Range("J132").Value = 3718.0908
$ws.Range("K132").Value = 14236.857
$ws.Range("L132").Value = 11154.2724
$ws.Range("M132").Value = -11706.857
$ws.Range("N132").Value = -16214.2724

$ws.Range("H136").Value = 4545.316
$ws.Range("I136").Value = 2040.8889
$ws.Range("K136").Value = 6122.6667
$ws.Range("M136").Value = -3572.6667

$ws.Range("H138").Value = 77679.75
$ws.Range("J138").Value = 77679.75
$ws.Range("L138").Value = 77679.75
$ws.Range("N138").Value = -87959.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 44098.668
$ws.Range("I51").Value = 35314.832
$ws.Range("K51").Value = 35314.832
$ws.Range("M51").Value = -34804.832

$ws.Range("H52").Value = 28995
$ws.Range("I52").Value = 28990
$ws.Range("K52").Value = 28990
$ws.Range("M52").Value = -28764

$ws.Range("H55").Value = 11112.25
$ws.Range("J55").Value = 11112.25
$ws.Range("L55").Value = 11112.25
$ws.Range("N55").Value = -11666.25

$ws.Range("H81").Value = 3931.2856
$ws.Range("I81").Value = 748.5
$ws.Range("K81").Value = 1497
$ws.Range("M81").Value = -436

$ws.Range("H84").Value = 3931.2856
$ws.Range("I84").Value = 748.5
$ws.Range("K84").Value = 7485
$ws.Range("M84").Value = -2181

$ws.Range("H95").Value = 15172
$ws.Range("J95").Value = 15172
$ws.Range("L95").Value = 15172
$ws.Range("N95").Value = -20664

$ws.Range("H113").Value = 1436.9231
$ws.Range("I113").Value = 906.9
$ws.Range("K113").Value = 2720.7
$ws.Range("M113").Value = -550.6999999999998

$ws.Range("H122").Value = 3656.9363
$ws.Range("I122").Value = 2868.8235
$ws.Range("K122").Value = 8606.470499999999
$ws.Range("M122").Value = -6156.470499999999

$ws.Range("H123").Value = 57249.5
$ws.Range("J123").Value = 57249.5
$ws.Range("L123").Value = 57249.5
$ws.Range("N123").Value = -67049.5

$ws.Range("H125").Value = 25357.5
$ws.Range("J125").Value = 25357.5
$ws.Range("L125").Value = 25357.5
$ws.Range("N125").Value = -35197.5

$ws.Range("H126").Value = 4086.842
$ws.Range("I126").Value = 3060.5386
$ws.Range("K126").Value = 9181.6158
$ws.Range("M126").Value = -6711.6158

$ws.Range("H132").Value = 3478.9756
$ws.Range("I132").Value = 1537.1482
$ws.Range("J132").Value = 7223.9287
$ws.Range("K132").Value = 4611.444600000001
$ws.Range("L132").Value = 21671.7861
$ws.Range("M132").Value = -2081.444600000001
$ws.Range("N132").Value = -26731.7861

$ws.Range("H136").Value = 9284.281000000001
$ws.Range("I136").Value = 2437.3333
$ws.Range("J136").Value = 12160
$ws.Range("K136").Value = 7311.999899999999
$ws.Range("L136").Value = 36480
$ws.Range("M136").Value = -4761.999899999999
$ws.Range("N136").Value = -41580

Write-Host "Applied scheduled runner profit updates across 8 sheets (81 rows)."
